$d = $word.ActiveDocument

$d.Content.Find.Execute("330×8=2640", $true, $false, $false, $false, $false, $true, 1, $false, "955×5=4775", 2) | Out-Null
$d.Content.Find.Execute("282×8=2256", $true, $false, $false, $false, $false, $true, 1, $false, "451×9=4059", 2) | Out-Null
$d.Content.Find.Execute("117×9=1053", $true, $false, $false, $false, $false, $true, 1, $false, "121×6=726", 2) | Out-Null
$d.Content.Find.Execute("145×3=435", $true, $false, $false, $false, $false, $true, 1, $false, "566×3=1698", 2) | Out-Null
$d.Content.Find.Execute("473×9=4257", $true, $false, $false, $false, $false, $true, 1, $false, "375×4=1500", 2) | Out-Null
$d.Content.Find.Execute("772×5=3860", $true, $false, $false, $false, $false, $true, 1, $false, "653×3=1959", 2) | Out-Null
$d.Content.Find.Execute("975×3=2925", $true, $false, $false, $false, $false, $true, 1, $false, "650×4=2600", 2) | Out-Null
$d.Content.Find.Execute("542×9=4878", $true, $false, $false, $false, $false, $true, 1, $false, "927×8=7416", 2) | Out-Null
$d.Content.Find.Execute("821×9=7389", $true, $false, $false, $false, $false, $true, 1, $false, "231×3=693", 2) | Out-Null
$d.Content.Find.Execute("947×8=7576", $true, $false, $false, $false, $false, $true, 1, $false, "496×9=4464", 2) | Out-Null
$d.Content.Find.Execute("570×6=3420", $true, $false, $false, $false, $false, $true, 1, $false, "544×2=1088", 2) | Out-Null
$d.Content.Find.Execute("974×7=6818", $true, $false, $false, $false, $false, $true, 1, $false, "857×8=6856", 2) | Out-Null
$d.Content.Find.Execute("913×7=6391", $true, $false, $false, $false, $false, $true, 1, $false, "864×2=1728", 2) | Out-Null
$d.Content.Find.Execute("454×2=908", $true, $false, $false, $false, $false, $true, 1, $false, "390×9=3510", 2) | Out-Null
$d.Content.Find.Execute("986×9=8874", $true, $false, $false, $false, $false, $true, 1, $false, "743×4=2972", 2) | Out-Null
$d.Content.Find.Execute("151×5=755", $true, $false, $false, $false, $false, $true, 1, $false, "613×5=3065", 2) | Out-Null
$d.Content.Find.Execute("312×6=1872", $true, $false, $false, $false, $false, $true, 1, $false, "854×9=7686", 2) | Out-Null
$d.Content.Find.Execute("331×7=2317", $true, $false, $false, $false, $false, $true, 1, $false, "182×6=1092", 2) | Out-Null
$d.Content.Find.Execute("873×2=1746", $true, $false, $false, $false, $false, $true, 1, $false, "383×4=1532", 2) | Out-Null
$d.Content.Find.Execute("109×7=763", $true, $false, $false, $false, $false, $true, 1, $false, "980×6=5880", 2) | Out-Null
$d.Content.Find.Execute("837×6=5022", $true, $false, $false, $false, $false, $true, 1, $false, "978×3=2934", 2) | Out-Null
$d.Content.Find.Execute("516×6=3096", $true, $false, $false, $false, $false, $true, 1, $false, "993×6=5958", 2) | Out-Null
$d.Content.Find.Execute("628×4=2512", $true, $false, $false, $false, $false, $true, 1, $false, "733×4=2932", 2) | Out-Null
$d.Content.Find.Execute("106×9=954", $true, $false, $false, $false, $false, $true, 1, $false, "541×6=3246", 2) | Out-Null
$d.Content.Find.Execute("666×3=1998", $true, $false, $false, $false, $false, $true, 1, $false, "400×8=3200", 2) | Out-Null
